# Breadcrumbs.xlsx - "updated main GSC export data"
#
# The "Chart" sheet holds a rolling window of daily GSC export rows:
#   column A = date label (stored as literal TEXT, e.g. "2025-10-09")
#   column B = Invalid count
#   column C = Valid/Items count
#
# A fresh export advances the window by one day: the oldest date row is
# dropped and a new row for the newest date is appended, so every
# existing row's data shifts up by one position and a brand new row is
# added at the bottom with a new date + new Items value.
#
# NOTE: Writing a date-looking literal straight into `.Value` makes Excel
# "helpfully" reinterpret it as a real date serial (changing both the
# stored value and the cell's number-format style). To keep column A as
# plain text - exactly like the source file - we move values between
# cells with Copy/PasteSpecial(xlPasteValues), which carries the literal
# text through untouched and leaves styles alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$headerRow  = 1
$firstRow   = $headerRow + 1
$lastRow    = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1
$dateCol    = 1
$itemsCol   = 3

$newDate  = "2026-01-06"
$newItems = 27

$xlPasteValues = -4163

# Shift rows firstRow..lastRow-1 up by one: row r becomes what row r+1 was.
for ($r = $firstRow; $r -le ($lastRow - 1); $r++) {
    $ws.Cells.Item($r + 1, $dateCol).Copy()
    $ws.Cells.Item($r, $dateCol).PasteSpecial($xlPasteValues)

    $ws.Cells.Item($r + 1, $itemsCol).Copy()
    $ws.Cells.Item($r, $itemsCol).PasteSpecial($xlPasteValues)
}

# New last row gets the newest date (kept as literal text via the
# formula -> copy -> paste-values-in-place trick) and its Items count.
$newDateCell = $ws.Cells.Item($lastRow, $dateCol)
$newDateCell.Formula = "=""" + $newDate + """"
$newDateCell.Copy()
$newDateCell.PasteSpecial($xlPasteValues)

$ws.Cells.Item($lastRow, $itemsCol).Value = $newItems

$excel.CutCopyMode = 0
